$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 382, shifting existing rows 382:410 down to 383:411.
$ws.Rows.Item(382).Insert()

# Populate the newly inserted row 382 with the new record.
$ws.Cells.Item(382, 1).Value = 4
$ws.Cells.Item(382, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(382, 3).Value = "Los Lagos"
$ws.Cells.Item(382, 4).Value = 45013
$ws.Cells.Item(382, 5).Value = 10
$ws.Cells.Item(382, 6).Value = 100112043
$ws.Cells.Item(382, 7).Value = "Pepino ensalada"
$ws.Cells.Item(382, 8).Value = "Sin especificar"
$ws.Cells.Item(382, 9).Value = "Primera"
$ws.Cells.Item(382, 10).Value = 400
$ws.Cells.Item(382, 11).Value = 13000
$ws.Cells.Item(382, 12).Value = 13000
$ws.Cells.Item(382, 13).Value = 13000
$ws.Cells.Item(382, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(382, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(382, 16).Value = 217
$ws.Cells.Item(382, 17).Value = 60
$ws.Cells.Item(382, 18).Value = "Hortaliza"
